$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new switches ("n1-g620-005-vc67-f1" / "n1-g620-006-vc67-f2") were added to the
# NORD VC67 fabric block, right after the existing rows for that fabric (row 13) and
# before the NORD Cloud152 block (previously starting at row 14). Insert two blank
# rows there; this shifts all subsequent rows down by two, matching the diff exactly.
$ws.Rows("14:15").Insert()

# Fill in the data for the newly inserted rows.
$ws.Range("A14").Value = "NORD VC67"
$ws.Range("B14").Value = "n1-g620-005-vc67-f1"
$ws.Range("C14").Value = "10:00:88:94:71:60:75:63"
$ws.Range("D14").Value = "Brocade G620"
$ws.Range("E14").Value = "3H1J110 un37"

$ws.Range("A15").Value = "NORD VC67"
$ws.Range("B15").Value = "n1-g620-006-vc67-f2"
$ws.Range("C15").Value = "10:00:88:94:71:ce:bf:fa"
$ws.Range("D15").Value = "Brocade G620"
$ws.Range("E15").Value = "3H1J090 un37"

# Match the cell formatting used in the target file: columns C and E of the two new
# rows carry the same (alternate) style as column B, rather than the plain style used
# by column A/D. Copy that formatting over.
$ws.Range("B14").Copy()
$ws.Range("C14:C15").PasteSpecial(-4122)
$ws.Range("E14:E15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The author's last selection in the saved file was cell C22.
$null = $ws.Range("C22").Select()
